$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data for the 14 new rows (rows 24-37), columns A..I
$rows = @(
    @('BMI','Brandmeldinstallatie','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','Brandmeldinstallatie, B1063, T2-01-R2011'),
    @('BMI','AspiratiemelderASD','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','AspiratiemelderASD, B1063, T2-01-R2011'),
    @('BMI','VideoSmokeDetectie','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','VideoSmokeDetectie, B1063, T2-01-R2011'),
    @('BMI','LineaireOptischeDetectie','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','LineaireOptischeDetectie, B1063, T2-01-R2011'),
    @('BMI','LineaireThermischeDetectie','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','LineaireThermischeDetectie, B1063, T2-01-R2011'),
    @('BMI','Vlammenmelder','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','Vlammenmelder, B1063, T2-01-R2011'),
    @('BMI','Handbrandmelder','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','Handbrandmelder, B1063, T2-01-R2011'),
    @('BMI','AkoestischeEnOptischeSignaalgevers','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','AkoestischeEnOptischeSignaalgevers, B1063, T2-01-R2011'),
    @('BMI','ThermischeMelder','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','ThermischeMelder, B1063, T2-01-R2011'),
    @('BMI','MultisensorPuntmelder','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','MultisensorPuntmelder, B1063, T2-01-R2011'),
    @('BMI','IOModule','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','IOModule, B1063, T2-01-R2011'),
    @('BMI','ExterneEnergievoorziening-BrandmeldinstallatieVeldcomponent','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','ExterneEnergievoorziening-BrandmeldinstallatieVeldcomponent, B1063, T2-01-R2011'),
    @('BMI','NevenIndicator','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','NevenIndicator, B1063, T2-01-R2011'),
    @('BMI','NevenpaneelBrandmeldcentrale','Ja','^B\d{1,4}$','B[Brandmeldcentralennr]','B1063','^[A-Za-z0-9 -]+,[[aascode]],.*','[IRI], [AasCode], [Ruimtenummer]','NevenpaneelBrandmeldcentrale, B1063, T2-01-R2011')
)

$startRow = 24
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}

# Resize the table (ListObject) to cover the newly added rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I37"))

# Make "Conventies" the active sheet/tab (was "Versie toetsingsregel")
$ws.Activate()
$ws.Range("I41").Select()
